# Insert a new weekly price record as row 10 in the dataset.
# All subsequent rows (old 10..103) shift down by one (to 11..104),
# which matches the unified diff (dimension A1:R103 -> A1:R104).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 10 (and everything below it) down by one row.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new market observation.
$ws.Range("A10").Value = 3
$ws.Range("B10").Value = "Femacal de La Calera"
$ws.Range("C10").Value = "Coquimbo"
$ws.Range("D10").Value = 45282
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 100112022
$ws.Range("G10").Value = "Arveja Verde"
$ws.Range("H10").Value = "Perfection"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 50
$ws.Range("K10").Value = 33000
$ws.Range("L10").Value = 35000
$ws.Range("M10").Value = 34400
$ws.Range("N10").Value = "`$/saco 25 kilos"
$ws.Range("O10").Value = "Provincia de Quillota"
$ws.Range("P10").Value = 1376
$ws.Range("Q10").Value = 25
$ws.Range("R10").Value = "Hortaliza"
